$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "Punendra"
$ws.Range("B7").Value = 2200039115
$ws.Range("C7").Value = "2200039115@kluniversity.in"
$ws.Range("D7").Value = "Category-1"
$ws.Range("E7").Value = "Y22"

$ws.Hyperlinks.Add($ws.Range("C7"), "mailto:2200039115@kluniversity.in") | Out-Null

$ws.Range("C4").Copy() | Out-Null
$ws.Range("C7").PasteSpecial(-4122) | Out-Null

$ws.Range("E7").Select() | Out-Null
